# Update the tractor-data sheet: "Trakstar Tractors" -> "Swaraj Tractors",
# re-verified models/images, and 4 additional model rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Cells.Item(1, 1).Value = "Brand"
$ws.Cells.Item(1, 2).Value = "Model"
$ws.Cells.Item(1, 3).Value = "Images"

$brand = "Swaraj Tractors"

# Each entry: Model name, Images list (as text, matching source data format)
$rows = @(
    @("855 FE",              "['855FEimg0-swaraj-855-fe-1694259363.png', '855FEimg1-855-fe-1631015724.png', '855FEimg2-upload-1631015724-0.png']"),
    @("744 FE",              "['744FEimg0-swaraj-744-fe-1694259976.png', '744FEimg1-744-fe-1631014732.png', '744FEimg2-upload-1631014732-0.png']"),
    @("735 FE",              "['735FEimg0-735-fe-1631013770.png', '735FEimg1-mqdefault.png', '735FEimg2-mqdefault.png']"),
    @("Code",                "['Codeimg0-code-1646204065.png', 'Codeimg1-code-1646204065.png', 'Codeimg2-mqdefault.png']"),
    @("744 XT",              "['744XTimg0-swaraj-744-xt-1694260212.png', '744XTimg1-744-xt-1631015130.png', '744XTimg2-upload-1631015130-0.png']"),
    @("742 XT",              "['742XTimg0-swaraj-742-xt-1694259930.png', '742XTimg1-742-xt-1631014604.png', '742XTimg2-swaraj-742-xt-16942599300.png']"),
    @("735 FE E",            "['735FEEimg0-swaraj-735-fe-e-1701335706.png', '735FEEimg1-swaraj-735-fe-e-17013357060.png', '735FEEimg2-swaraj-735-fe-e-1701335706.png']"),
    @("969 FE Trem IV-4wd",  "['969FETremIV-4wdimg0-swaraj-969-fe-trem-iv-4wd-1689145707.png', '969FETremIV-4wdimg1-swaraj-969-fe-trem-iv-4wd-1689145707.png', '969FETremIV-4wdimg2-mqdefault.png']"),
    @("855 FE 4WD",          "['855FE4WDimg0-swaraj-855-fe-4wd-1695808511.png', '855FE4WDimg1-swaraj-855-fe-4wd-16958103590.png', '855FE4WDimg2-swaraj-855-fe-4wd-16958103590.png']"),
    @("744 FE 4WD",          "['744FE4WDimg0-744-fe-4wd-1631014809.png', '744FE4WDimg1-mqdefault.png', '744FE4WDimg2-mqdefault.png']"),
    @("735 XT",              "['735XTimg0-735-xt-1631014365.png', '735XTimg1-735-xt-1631014365.png', '735XTimg2-mqdefault.png']")
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $brand
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $r = $r + 1
}
